$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Force the cell to store the value as text (matching the
    # original inline-string cell type) instead of letting Excel
    # auto-convert number-looking strings (e.g. "572.22") into
    # a numeric cell. The leading apostrophe is the classic Excel
    # "treat as text" prefix; re-applying the Normal style afterwards
    # clears the quote-prefix formatting flag so no stray style
    # attribute gets attached to the cell.
    $cell = $ws.Range($cellAddr)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" '69.439.26'
$ws.Range("E2").Value = '  -1.44%  '
Set-TextValue "D3" '2.512.56'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue "D5" '572.22'
Set-TextValue "D6" '165.86'
$ws.Range("E6").Value = '  -2.14%  '
$ws.Range("E7").Value = '  -0.07%  '
Set-TextValue "D9" '2.510.34'
$ws.Range("E9").Value = '  -0.23%  '
$ws.Range("E10").Value = '  -2.42%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  +4.24%  '
$ws.Range("E13").Value = '  +1.06%  '
Set-TextValue "D14" '2.974.17'
$ws.Range("E14").Value = '  -0.24%  '
Set-TextValue "D15" '69.323.36'
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("E16").Value = '  -2.31%  '
Set-TextValue "D17" '24.78'
$ws.Range("E17").Value = '  -0.94%  '
Set-TextValue "D18" '2.517.71'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  -1.48%  '
Set-TextValue "D20" '7.55'
$ws.Range("E20").Value = '  -1.41%  '
Set-TextValue "D21" '348.61'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("E22").Value = '  -1.32%  '
$ws.Range("E23").Value = '  +0.44%  '
Set-TextValue "D25" '70.27'
$ws.Range("E25").Value = '  +1.54%  '
$ws.Range("E26").Value = '  -4.06%  '
Set-TextValue "D27" '8.86'
$ws.Range("E27").Value = '  -3.09%  '
Set-TextValue "D28" '2.658.57'
$ws.Range("E28").Value = '  +0.40%  '
Set-TextValue "D29" '0.995'
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  -3.29%  '
Set-TextValue "D31" '7.83'
$ws.Range("E31").Value = '  -0.27%  '
Set-TextValue "D32" '460.83'
$ws.Range("E32").Value = '  -4.67%  '
$ws.Range("E34").Value = '  -2.41%  '
Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +1.77%  '
Set-TextValue "D37" '157.25'
$ws.Range("E37").Value = '  +0.21%  '
Set-TextValue "D38" '19.07'
$ws.Range("E38").Value = '  +0.92%  '
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("E43").Value = '  -2.75%  '
Set-TextValue "D44" '38.14'
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("E45").Value = '  -7.37%  '
$ws.Range("E46").Value = '  -6.90%  '
Set-TextValue "D47" '141.76'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("E48").Value = '  -1.91%  '
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  -3.36%  '
